# "Modificada animacao de alteracao de atributo por setter"
#
# Re-lays-out the "setEndereco" attribute/setter diagram on slide 1:
# several shapes are nudged and resized, three labels get new font
# sizes, and the "cliente1" label is enlarged & retyped.
#
# NOTE: two hunks of the source diff are not reachable from the
# PowerPoint COM object model and are intentionally skipped:
#   * presentation.xml: an *empty* <p15:sldGuideLst/> ext block - this
#     is written by the native app as a side effect of internal guide
#     bookkeeping; Presentation.Guides / Master.Guides are stubbed out
#     (return Nothing) in this host and there is no Add()/Count() path
#     that serialises anything.
#   * the <a:endParaRPr>/smtClean="0" artifacts that real PowerPoint
#     leaves behind after an in-place retype: TextRange.Font only ever
#     touches the run's <a:rPr>, there is no addressable handle onto a
#     paragraph's trailing end-of-paragraph run.
# Everything else (every shape's position/size and the four run font
# sizes) is applied below.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# Rectangle 3 - taller gap at the top, shorter box overall
# ---------------------------------------------------------------
$rect3 = $s.Shapes.Item("Rectangle 3")
$rect3.Top = 29.037639617919922
$rect3.Height = 166.73143005371094

# ---------------------------------------------------------------
# Group 13 ("Rectangle 4" + "TextBox 5") - the inner "endereco" box
# shrinks, so the group's own bounding box shrinks with it.
# ---------------------------------------------------------------
$group13 = $s.Shapes.Item("Group 13")
$textBox5 = $group13.GroupItems.Item("TextBox 5")
$textBox5.Height = 47.96866226196289
$textBox5.TextFrame.TextRange.Font.Size = 12

# Keep the group wrapper's own size in sync with the resized child
# (this host does not auto-recompute group bounds from member edits).
$group13.Height = 40.45716857910156

# ---------------------------------------------------------------
# TextBox 12 ("cliente1") - enlarged, re-centred, retyped at 20pt bold
# ---------------------------------------------------------------
$textBox12 = $s.Shapes.Item("TextBox 12")
$textBox12.Left = 28.745119094848633
$textBox12.Top = -1.7053543329238892
$textBox12.Width = 80.88724517822266
$textBox12.Height = 31.504724502563477
$textBox12.TextFrame.TextRange.Font.Size = 20
$textBox12.TextFrame.TextRange.Font.Bold = $true

# ---------------------------------------------------------------
# TextBox 17 ("setEndereco") - moved/enlarged, font bumped to 16pt
# ---------------------------------------------------------------
$textBox17 = $s.Shapes.Item("TextBox 17")
$textBox17.Left = 22.707717895507812
$textBox17.Top = 101.67449188232422
$textBox17.Width = 96.6768569946289
$textBox17.Height = 26.65779685974121
$textBox17.TextFrame.TextRange.Font.Size = 16

# ---------------------------------------------------------------
# Rectangle 27 - shifted up
# ---------------------------------------------------------------
$rect27 = $s.Shapes.Item("Rectangle 27")
$rect27.Top = 134.34095764160156

# ---------------------------------------------------------------
# TextBox 28 ("endereco") - shifted up, slightly shorter, font 12pt
# ---------------------------------------------------------------
$textBox28 = $s.Shapes.Item("TextBox 28")
$textBox28.Top = 152.37158203125
$textBox28.Height = 21.810945510864258
$textBox28.TextFrame.TextRange.Font.Size = 12

# ---------------------------------------------------------------
# Double Bracket 2 - shifted up
# ---------------------------------------------------------------
$bracket2 = $s.Shapes.Item("Double Bracket 2")
$bracket2.Top = 116.14803314208984

# ---------------------------------------------------------------
# TextBox 1 ("Manoel", lower copy) - shifted up
# ---------------------------------------------------------------
$textBox1 = $s.Shapes.Item("TextBox 1")
$textBox1.Top = 135.3643341064453

# ---------------------------------------------------------------
# Curved Up Arrow 23 - slides right/up and shortens along its length
# ---------------------------------------------------------------
$arrow23 = $s.Shapes.Item("Curved Up Arrow 23")
$arrow23.Left = 67.64968872070312
$arrow23.Top = 87.84795379638672
$arrow23.Width = 103.27228546142578
